# Commit: "Natmi following Dr Hou advice"
# Update NATMI LR-pair stats for Tnc-Egfr sheet: ligand/receptor-expressing
# cell counts (E, K) change from 1 to 3 for every data row, which cascades
# into the dependent expression/specificity/edge-weight columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per data row (2-10) for the columns that change.
# Columns A-D (labels), F, and L are unchanged by this edit.
$updates = @{
    2 = @{ "E" = 3; "G" = 12.07002266666666; "H" = 36.21006799999999; "I" = 0.7601982364861632; "J" = 0.7601982364861634; "K" = 3; "M" = 2.08532; "N" = 6.25596; "O" = 0.01753772176136817; "P" = 0.01753772176136816; "Q" = 25.16985966725333; "R" = 226.52873700528; "S" = 0.01333214515497709; "T" = 0.01333214515497709 }
    3 = @{ "E" = 3; "G" = 12.07002266666666; "H" = 36.21006799999999; "I" = 0.7601982364861632; "J" = 0.7601982364861634; "K" = 3; "M" = 101.898173; "N" = 305.694519; "O" = 0.8569724579756384; "P" = 0.8569724579756383; "Q" = 1229.913257801921; "R" = 11069.21932021729; "S" = 0.6514689512702929; "T" = 0.6514689512702929 }
    4 = @{ "E" = 3; "G" = 12.07002266666666; "H" = 36.21006799999999; "I" = 0.7601982364861632; "J" = 0.7601982364861634; "K" = 3; "M" = 14.921347; "N" = 44.76404100000001; "O" = 0.1254898202629935; "P" = 0.1254898202629935; "Q" = 180.1009965071987; "R" = 1620.908968564788; "S" = 0.09539714006089327; "T" = 0.09539714006089325 }
    5 = @{ "E" = 3; "G" = 1.308268; "H" = 3.924804; "I" = 0.08239777620284613; "J" = 0.08239777620284613; "K" = 3; "M" = 2.08532; "N" = 6.25596; "O" = 0.01753772176136817; "P" = 0.01753772176136816; "Q" = 2.72815742576; "R" = 24.55341683184; "S" = 0.001445069272800999; "T" = 0.001445069272800998 }
    6 = @{ "E" = 3; "G" = 1.308268; "H" = 3.924804; "I" = 0.08239777620284613; "J" = 0.08239777620284613; "K" = 3; "M" = 101.898173; "N" = 305.694519; "O" = 0.8569724579756384; "P" = 0.8569724579756383; "Q" = 133.310118994364; "R" = 1199.791070949276; "S" = 0.07061262480427961; "T" = 0.0706126248042796 }
    7 = @{ "E" = 3; "G" = 1.308268; "H" = 3.924804; "I" = 0.08239777620284613; "J" = 0.08239777620284613; "K" = 3; "M" = 14.921347; "N" = 44.76404100000001; "O" = 0.1254898202629935; "P" = 0.1254898202629935; "Q" = 19.521120796996; "R" = 175.690087172964; "S" = 0.01034008212576553; "T" = 0.01034008212576552 }
    8 = @{ "E" = 3; "G" = 2.499176666666667; "H" = 7.49753; "I" = 0.1574039873109905; "J" = 0.1574039873109906; "K" = 3; "M" = 2.08532; "N" = 6.25596; "O" = 0.01753772176136817; "P" = 0.01753772176136816; "Q" = 5.211583086533333; "R" = 46.9042477788; "S" = 0.002760507333590077; "T" = 0.002760507333590077 }
    9 = @{ "E" = 3; "G" = 2.499176666666667; "H" = 7.49753; "I" = 0.1574039873109905; "J" = 0.1574039873109906; "K" = 3; "M" = 101.898173; "N" = 305.694519; "O" = 0.8569724579756384; "P" = 0.8569724579756383; "Q" = 254.6615363375633; "R" = 2291.95382703807; "S" = 0.1348908819010657; "T" = 0.1348908819010658 }
    10 = @{ "E" = 3; "G" = 2.499176666666667; "H" = 7.49753; "I" = 0.1574039873109905; "J" = 0.1574039873109906; "K" = 3; "M" = 14.921347; "N" = 44.76404100000001; "O" = 0.1254898202629935; "P" = 0.1254898202629935; "Q" = 37.29108225763667; "R" = 335.6197403187301; "S" = 0.01975259807633472; "T" = 0.01975259807633471 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($colLetter in $rowData.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowData[$colLetter]
    }
}
